# Correct CTDMO reference designators: rows 65-68 and 70-73 had duplicated
# designators (CTDMOG044 / CTDMOG045 were entered twice); they should read
# CTDMOG046 and CTDMOG047 respectively.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Asset_Cal_Info")

# --- Step 1: Fix the duplicate reference designators (the core data correction) ---
$ws.Range("A65").Value2 = "GA03FLMA-RIM01-02-CTDMOG046"
$ws.Range("A66").Value2 = "GA03FLMA-RIM01-02-CTDMOG046"
$ws.Range("A67").Value2 = "GA03FLMA-RIM01-02-CTDMOG046"
$ws.Range("A68").Value2 = "GA03FLMA-RIM01-02-CTDMOG046"
$ws.Range("A70").Value2 = "GA03FLMA-RIM01-02-CTDMOG047"
$ws.Range("A71").Value2 = "GA03FLMA-RIM01-02-CTDMOG047"
$ws.Range("A72").Value2 = "GA03FLMA-RIM01-02-CTDMOG047"
$ws.Range("A73").Value2 = "GA03FLMA-RIM01-02-CTDMOG047"

# --- Step 2: Build the helper "M" column (unique reference designator list) ---
$ws.Range("M2").Value2 = "GA03FLMA-RIS01-00-SIOENG000"
$ws.Range("M3").Value2 = "GA03FLMA-RIS01-05-FLORTD000"
$ws.Range("M4").Value2 = "GA03FLMA-RIS01-04-PHSENF000"
$ws.Range("M5").Value2 = "GA03FLMA-RIS01-03-DOSTAD000"
$ws.Range("M6").Value2 = "GA03FLMA-RIM01-02-ADCPSL003"
$ws.Range("M7").Value2 = "GA03FLMA-RIM01-02-CTDMOG040"
$ws.Range("M8").Value2 = "GA03FLMA-RIM01-02-CTDMOG041"
$ws.Range("M9").Value2 = "GA03FLMA-RIM01-02-CTDMOG042"
$ws.Range("M10").Value2 = "GA03FLMA-RIM01-02-CTDMOG043"
$ws.Range("M11").Value2 = "GA03FLMA-RIM01-02-CTDMOG044"
$ws.Range("M12").Value2 = "GA03FLMA-RIM01-02-CTDMOG045"
$ws.Range("M13").Value2 = "GA03FLMA-RIM01-02-CTDMOG046"
$ws.Range("M14").Value2 = "GA03FLMA-RIM01-02-CTDMOG047"
$ws.Range("M15").Value2 = "GA03FLMA-RIM01-02-CTDMOG048"
$ws.Range("M16").Value2 = "GA03FLMA-RIM01-02-CTDMOH049"
$ws.Range("M17").Value2 = "GA03FLMA-RIM01-02-CTDMOH050"
$ws.Range("M18").Value2 = "GA03FLMA-RIM01-02-CTDMOH051"
$ws.Range("M19").Value2 = "GA03FLMA-RIM01-00-SIOENG000"

# --- Step 3: Build the helper "L" column (MATCH of each row against M:M) ---
$ws.Range("L2").Formula = "=MATCH(A2,M:M,0)"
$ws.Range("L3").Formula = "=MATCH(A3,M:M,0)"
$ws.Range("L4").Formula = "=MATCH(A4,M:M,0)"
$ws.Range("L5").Formula = "=MATCH(A5,M:M,0)"
$ws.Range("L6").Formula = "=MATCH(A6,M:M,0)"
$ws.Range("L7").Formula = "=MATCH(A7,M:M,0)"
$ws.Range("L8").Formula = "=MATCH(A8,M:M,0)"
$ws.Range("L9").Formula = "=MATCH(A9,M:M,0)"
$ws.Range("L10").Formula = "=MATCH(A10,M:M,0)"
$ws.Range("L11").Formula = "=MATCH(A11,M:M,0)"
$ws.Range("L12").Formula = "=MATCH(A12,M:M,0)"
$ws.Range("L13").Formula = "=MATCH(A13,M:M,0)"
$ws.Range("L14").Formula = "=MATCH(A14,M:M,0)"
$ws.Range("L15").Formula = "=MATCH(A15,M:M,0)"
$ws.Range("L16").Formula = "=MATCH(A16,M:M,0)"
$ws.Range("L17").Formula = "=MATCH(A17,M:M,0)"
$ws.Range("L18").Formula = "=MATCH(A18,M:M,0)"
$ws.Range("L19").Formula = "=MATCH(A19,M:M,0)"
$ws.Range("L20").Formula = "=MATCH(A20,M:M,0)"
$ws.Range("L21").Formula = "=MATCH(A21,M:M,0)"
$ws.Range("L22").Formula = "=MATCH(A22,M:M,0)"
$ws.Range("L23").Formula = "=MATCH(A23,M:M,0)"
$ws.Range("L24").Formula = "=MATCH(A24,M:M,0)"
$ws.Range("L25").Formula = "=MATCH(A25,M:M,0)"
$ws.Range("L26").Formula = "=MATCH(A26,M:M,0)"
$ws.Range("L27").Formula = "=MATCH(A27,M:M,0)"
$ws.Range("L28").Formula = "=MATCH(A28,M:M,0)"
$ws.Range("L29").Formula = "=MATCH(A29,M:M,0)"
$ws.Range("L30").Formula = "=MATCH(A30,M:M,0)"
$ws.Range("L31").Formula = "=MATCH(A31,M:M,0)"
$ws.Range("L32").Formula = "=MATCH(A32,M:M,0)"
$ws.Range("L33").Formula = "=MATCH(A33,M:M,0)"
$ws.Range("L34").Formula = "=MATCH(A34,M:M,0)"
$ws.Range("L35").Formula = "=MATCH(A35,M:M,0)"
$ws.Range("L36").Formula = "=MATCH(A36,M:M,0)"
$ws.Range("L37").Formula = "=MATCH(A37,M:M,0)"
$ws.Range("L38").Formula = "=MATCH(A38,M:M,0)"
$ws.Range("L39").Formula = "=MATCH(A39,M:M,0)"
$ws.Range("L40").Formula = "=MATCH(A40,M:M,0)"
$ws.Range("L41").Formula = "=MATCH(A41,M:M,0)"
$ws.Range("L42").Formula = "=MATCH(A42,M:M,0)"
$ws.Range("L43").Formula = "=MATCH(A43,M:M,0)"
$ws.Range("L44").Formula = "=MATCH(A44,M:M,0)"
$ws.Range("L45").Formula = "=MATCH(A45,M:M,0)"
$ws.Range("L46").Formula = "=MATCH(A46,M:M,0)"
$ws.Range("L47").Formula = "=MATCH(A47,M:M,0)"
$ws.Range("L48").Formula = "=MATCH(A48,M:M,0)"
$ws.Range("L49").Formula = "=MATCH(A49,M:M,0)"
$ws.Range("L50").Formula = "=MATCH(A50,M:M,0)"
$ws.Range("L51").Formula = "=MATCH(A51,M:M,0)"
$ws.Range("L52").Formula = "=MATCH(A52,M:M,0)"
$ws.Range("L53").Formula = "=MATCH(A53,M:M,0)"
$ws.Range("L54").Formula = "=MATCH(A54,M:M,0)"
$ws.Range("L55").Formula = "=MATCH(A55,M:M,0)"
$ws.Range("L56").Formula = "=MATCH(A56,M:M,0)"
$ws.Range("L57").Formula = "=MATCH(A57,M:M,0)"
$ws.Range("L58").Formula = "=MATCH(A58,M:M,0)"
$ws.Range("L59").Formula = "=MATCH(A59,M:M,0)"
$ws.Range("L60").Formula = "=MATCH(A60,M:M,0)"
$ws.Range("L61").Formula = "=MATCH(A61,M:M,0)"
$ws.Range("L62").Formula = "=MATCH(A62,M:M,0)"
$ws.Range("L63").Formula = "=MATCH(A63,M:M,0)"
$ws.Range("L64").Formula = "=MATCH(A64,M:M,0)"
$ws.Range("L65").Formula = "=MATCH(A65,M:M,0)"
$ws.Range("L66").Formula = "=MATCH(A66,M:M,0)"
$ws.Range("L67").Formula = "=MATCH(A67,M:M,0)"
$ws.Range("L68").Formula = "=MATCH(A68,M:M,0)"
$ws.Range("L69").Formula = "=MATCH(A69,M:M,0)"
$ws.Range("L70").Formula = "=MATCH(A70,M:M,0)"
$ws.Range("L71").Formula = "=MATCH(A71,M:M,0)"
$ws.Range("L72").Formula = "=MATCH(A72,M:M,0)"
$ws.Range("L73").Formula = "=MATCH(A73,M:M,0)"
$ws.Range("L74").Formula = "=MATCH(A74,M:M,0)"
$ws.Range("L75").Formula = "=MATCH(A75,M:M,0)"
$ws.Range("L76").Formula = "=MATCH(A76,M:M,0)"
$ws.Range("L77").Formula = "=MATCH(A77,M:M,0)"
$ws.Range("L78").Formula = "=MATCH(A78,M:M,0)"
$ws.Range("L79").Formula = "=MATCH(A79,M:M,0)"
$ws.Range("L80").Formula = "=MATCH(A80,M:M,0)"
$ws.Range("L81").Formula = "=MATCH(A81,M:M,0)"
$ws.Range("L82").Formula = "=MATCH(A82,M:M,0)"
$ws.Range("L83").Formula = "=MATCH(A83,M:M,0)"
$ws.Range("L84").Formula = "=MATCH(A84,M:M,0)"
$ws.Range("L85").Formula = "=MATCH(A85,M:M,0)"
$ws.Range("L86").Formula = "=MATCH(A86,M:M,0)"
$ws.Range("L87").Formula = "=MATCH(A87,M:M,0)"
$ws.Range("L88").Formula = "=MATCH(A88,M:M,0)"
$ws.Range("L89").Formula = "=MATCH(A89,M:M,0)"
$ws.Range("L90").Formula = "=MATCH(A90,M:M,0)"
$ws.Range("L91").Formula = "=MATCH(A91,M:M,0)"
$ws.Range("L92").Formula = "=MATCH(A92,M:M,0)"
$ws.Range("L93").Formula = "=MATCH(A93,M:M,0)"
$ws.Range("L94").Formula = "=MATCH(A94,M:M,0)"
$ws.Range("L95").Formula = "=MATCH(A95,M:M,0)"
$ws.Range("L96").Formula = "=MATCH(A96,M:M,0)"
$ws.Range("L97").Formula = "=MATCH(A97,M:M,0)"

# --- Step 4: Build the helper "Q" column (reverse MATCH of each M value against A:A) ---
$ws.Range("Q2").Formula = "=MATCH(M2,A:A,0)"
$ws.Range("Q3").Formula = "=MATCH(M3,A:A,0)"
$ws.Range("Q4").Formula = "=MATCH(M4,A:A,0)"
$ws.Range("Q5").Formula = "=MATCH(M5,A:A,0)"
$ws.Range("Q6").Formula = "=MATCH(M6,A:A,0)"
$ws.Range("Q7").Formula = "=MATCH(M7,A:A,0)"
$ws.Range("Q8").Formula = "=MATCH(M8,A:A,0)"
$ws.Range("Q9").Formula = "=MATCH(M9,A:A,0)"
$ws.Range("Q10").Formula = "=MATCH(M10,A:A,0)"
$ws.Range("Q11").Formula = "=MATCH(M11,A:A,0)"
$ws.Range("Q12").Formula = "=MATCH(M12,A:A,0)"
$ws.Range("Q13").Formula = "=MATCH(M13,A:A,0)"
$ws.Range("Q14").Formula = "=MATCH(M14,A:A,0)"
$ws.Range("Q15").Formula = "=MATCH(M15,A:A,0)"
$ws.Range("Q16").Formula = "=MATCH(M16,A:A,0)"
$ws.Range("Q17").Formula = "=MATCH(M17,A:A,0)"
$ws.Range("Q18").Formula = "=MATCH(M18,A:A,0)"
$ws.Range("Q19").Formula = "=MATCH(M19,A:A,0)"

# --- Step 5: row-height touch-up around the corrected rows ---
$ws.Rows.Item(73).RowHeight = 16.5
$ws.Rows.Item(74).RowHeight = 16.5

# --- Step 6: restore the view state (scroll position / active selection) ---
$ws.Application.ActiveWindow.ScrollRow = 58
$ws.Range("N77").Select()
